$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Sheet1" to "Table"
$ws.Name = "Table"

# Update the selection/active cell on the sheet
$ws.Activate()
$ws.Range("E12").Select()
